$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 285, pushing existing rows 285:303 down to 286:304
$ws.Rows("285:285").Insert()

# Populate the newly inserted row 285 with the new data record
$ws.Range("A285").Value = 3
$ws.Range("B285").Value = "Femacal de La Calera"
$ws.Range("C285").Value = "Coquimbo"
$ws.Range("D285").Value = 44610
$ws.Range("E285").Value = 5
$ws.Range("F285").Value = 100112040
$ws.Range("G285").Value = "Cilantro"
$ws.Range("H285").Value = "Sin especificar"
$ws.Range("I285").Value = "Primera"
$ws.Range("J285").Value = 105
$ws.Range("K285").Value = 4500
$ws.Range("L285").Value = 5000
$ws.Range("M285").Value = 4762
$ws.Range("N285").Value = "$/docena de atados (3 kilos)"
$ws.Range("O285").Value = "Provincia de Quillota"
$ws.Range("P285").Value = 1587
$ws.Range("Q285").Value = 3
$ws.Range("R285").Value = "Hortaliza"
